$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 19:07:23"
$wsZhCn.Range("H2").Value = "2016-03-21 19:07:44"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 19:07:26"
$wsDeDe.Range("H2").Value = "2016-03-21 19:07:50"
